$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ("253 Hoddle Street Apartment Complex Collingwood", 5),
    ("3535 Opal Meadow Heights Aged Care Community Meadow Heights", 26),
    ("A1 Bakery Brunswick", 5),
    ("Acquire BPO Southbank", 8),
    ("Al Haj Halal Meats Glenroy", 43),
    ("Al-Taqwa College Truganina", 7),
    ("Amiga Montessori Craigieburn", 25),
    ("Best&Less Fountain Gate Narre Warren", 7),
    ("Budget Car and Truck Rentals Campbellfield", 7),
    ("CS Square Caroline Springs", 7),
    ("Campbellfield Heights Primary School Campbellfield", 6),
    ("Cannie Road Construction Site Cannie", 9),
    ("Caroline Springs Police Station", 11),
    ("Cedars Medical Clinic Coburg", 44),
    ("Chemist Warehouse Fillo Drive Somerton", 5),
    ("City of Hobsons Bay Community", 8),
    ("City of Moreland Community", 7),
    ("City of Wyndham Community", 7),
    ("Classy Cabinets and Kitchens Craigieburn", 11),
    ("Coles Aurora Village Epping", 6),
    ("Coles Broadmeadows Central Shopping Centre", 10),
    ("Coles Campbellfield Plaza Campbellfield", 10),
    ("Coles Coburg North Village", 27),
    ("Coles Greenvale Shopping Centre", 7),
    ("Coles Pakenham Place Shopping Centre", 8),
    ("Coles Roxburgh Village Roxburgh Park", 11),
    ("Community Kids Meadow Heights", 16),
    ("Construction Site Olea Apartment Caulfield North", 10),
    ("Costco Wholesale Epping", 18),
    ("DRC Laverton Automotive Repairs Laverton North", 5),
    ("Direct Freight Express Campbellfield", 12),
    ("Elite Smart Community Care Campbellfield", 5),
    ("Fitzroy Community School Fitzroy North", 45),
    ("Fonterra Manufacturing Workplace Campbellfield", 5),
    ("Glenroy West Primary School", 5),
    ("Hamilton Marino 236 Jasper Road McKinnon", 11),
    ("Health Care Providers Association South Melbourne", 13),
    ("IGA Meadow Heights Shopping Centre Meadow Heights", 7),
    ("ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine", 7),
    ("Ilim College Glenroy Campus Hadfield", 22),
    ("Ilim Learning Sanctuary Glenroy", 15),
    ("Industrial Galvanizers Valmont Coatings Campbellfield", 17),
    ("Islamic College of Melbourne Tarneit", 8),
    ("KFC Fawkner", 8),
    ("Kasr Sweets Coolaroo", 7),
    ("Kids House Early Learning Cheltenham", 13),
    ("Learning Nest Early Learning Centre Meadow Heights", 8),
    ("Level Crossing Removal Project Lilydale Construction Site John Street", 8),
    ("Lineage Logistics Laverton North", 8),
    ("Linfox Somerton National Distribution Centre Somerton", 9),
    ("McDonalds Thomastown II", 7),
    ("Melbourne Metropolitan Remand Centre Ravenhall", 10),
    ("Melbourne Truck Repairs Campbellfield", 8),
    ("Melbourne West Police Station Docklands", 6),
    ("Melbourne Youth Justice Centre Parkville", 6),
    ("Melton Police Station Melton", 5),
    ("Mercy Hospital for Women Heidelberg", 5),
    ("Mernda YMCA Early Learning Centre Mernda", 5),
    ("Montessori Beginnings Greenvale", 6),
    ("MyCentre Childcare Broadmeadows", 15),
    ("National Gallery of Victoria Melbourne", 9),
    ("Newbury Child and Community Centre Craigieburn", 9),
    ("Nido Early School Moonee Ponds", 16),
    ("Nino Early Learning Adventures Lalor", 5),
    ("North Geelong House Party", 7),
    ("Northern Health Northern Hospital Epping Emergency Department Tier 1B", 46),
    ("OnQ Plumbing and Excavations Craigieburn", 18),
    ("Oporto Coolaroo", 15),
    ("Oscar Romero Catholic Primary School Craigieburn", 5),
    ("Paisley Park Early Learning Centre Bundoora", 10),
    ("Panorama Construction Site Whitehorse Rd Box Hill", 20),
    ("People First Healthcare Home Residence Disability Support Taylors Lakes", 5),
    ("Ramsay Health Care Warringal Private Hospital Heidelberg", 9),
    ("Richmond Quarter 261-271 Bridge Road Construction Site Richmond", 9),
    ("Salta Drive Construction Site Rangedale Drainage Altona North", 6),
    ("Serco Mill Park", 5),
    ("St Vincents Hospital Emergency Department Melbourne", 6),
    ("Tek Foods Somerton", 13),
    ("The Homestead Child and Family Centre Roxburgh Park", 17),
    ("The Royal Children's Hospital Melbourne Emergency Department Parkville Tier 1B", 9),
    ("ThorwestenCabinets Pakenham", 12),
    ("Tip Top Warehouse Dandenong", 5),
    ("Total Window Concepts Hoppers Crossing", 6),
    ("Tunstall Fresh Tunstall Square Shopping Centre Doncaster East", 5),
    ("Unilodge College Square Student Accommodation 570 Lygon Street Carlton", 14),
    ("Werribee Mercy Hospital Emergency Department", 7),
    ("Western Health Footscray Hospital Emergency Department", 5),
    ("Western Health Sunshine Hospital Emergency Department", 7),
    ("Woodlands Long Day Care and Kindergarten Roxburgh Park", 5),
    ("Woolworths Greenvale Lakes Roxburgh Park", 6),
    ("Yara Childcare Centre Truganina", 10)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

Write-Output "updated $($data.Count) data rows (now rows 2..$(1+$data.Count))"